$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Sending cluster: ECs, Target cluster: FAPs) ---
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 2.708836333333334
$ws.Range("H2").Value = 8.126509
$ws.Range("I2").Value = 0.0171826329450544
$ws.Range("J2").Value = 0.0171826329450544
$ws.Range("M2").Value = 0.001809666666666667
$ws.Range("N2").Value = 0.005429
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.004902090817888889
$ws.Range("R2").Value = 0.044118817361
$ws.Range("S2").Value = 0.0171826329450544
$ws.Range("T2").Value = 0.0171826329450544

# --- Row 3 (Sending cluster: FAPs, Target cluster: FAPs) ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 100.4067206666667
$ws.Range("H3").Value = 301.220162
$ws.Range("I3").Value = 0.6368977723762839
$ws.Range("J3").Value = 0.6368977723762839
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.1817026954997778
$ws.Range("R3").Value = 1.635324259498
$ws.Range("S3").Value = 0.6368977723762839
$ws.Range("T3").Value = 0.6368977723762839

# --- Row 4 (Sending cluster: MuSCs, Target cluster: FAPs) ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 54.53410833333334
$ws.Range("H4").Value = 163.602325
$ws.Range("I4").Value = 0.3459195946786617
$ws.Range("J4").Value = 0.3459195946786617
$ws.Range("M4").Value = 0.001809666666666667
$ws.Range("N4").Value = 0.005429
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.09868855804722222
$ws.Range("R4").Value = 0.888197022425
$ws.Range("S4").Value = 0.3459195946786617
$ws.Range("T4").Value = 0.3459195946786617

# --- Remove rows 5, 6, 7 (now obsolete data) ---
$ws.Rows("5:7").Delete()
